$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: FAPs / Cxcl5 / Cxcr1 / ECs
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Cxcl5"
$ws.Range("C2").Value = "Cxcr1"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 2.477721333333333
$ws.Range("H2").Value = 7.433164
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.0003833333333333333
$ws.Range("N2").Value = 0.00115
$ws.Range("O2").Value = 0.005317132262509131
$ws.Range("P2").Value = 0.005317132262509131
$ws.Range("Q2").Value = 0.0009497931777777777
$ws.Range("R2").Value = 0.0085481386
$ws.Range("S2").Value = 0.005317132262509131
$ws.Range("T2").Value = 0.005317132262509131

# Row 3: FAPs / Cxcl5 / Cxcr1 / Resolving-Mac
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Cxcl5"
$ws.Range("C3").Value = "Cxcr1"
$ws.Range("D3").Value = "Resolving-Mac"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 2.477721333333333
$ws.Range("H3").Value = 7.433164
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.07171066666666667
$ws.Range("N3").Value = 0.215132
$ws.Range("O3").Value = 0.9946828677374909
$ws.Range("P3").Value = 0.9946828677374908
$ws.Range("Q3").Value = 0.1776790486275556
$ws.Range("R3").Value = 1.599111437648
$ws.Range("S3").Value = 0.9946828677374909
$ws.Range("T3").Value = 0.9946828677374908

# Remove old rows 4 and 5 (data no longer present)
$ws.Range("A4:T5").Delete()
